$d = $word.ActiveDocument

# Locate the template placeholder text "{getFundingDocInfo()}" in the document body.
$hit = $d.Content
$found = $hit.Find.Execute("{getFundingDocInfo()}", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text '{getFundingDocInfo()}'"
}
$hitStart = $hit.Start
$hitEnd = $hit.End

# Re-seat a fresh Range over the exact hit so InsertXML replaces it in place.
$r = $d.Range($hitStart, $hitEnd)

# Common run formatting shared by every run in this span (bold, underlined Times New Roman 10pt).
$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/></w:rPr>'

$innerXml = (
    '<w:r>' + $rPr + '<w:t>{</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $rPr + '<w:t>f</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>unding</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>Doc</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>Inf</w:t></w:r>' +
    '<w:r>' + $rPr + '<w:t>o</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $rPr + '<w:t>}</w:t></w:r>'
)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' + $innerXml + '</w:p></w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
